$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1865
$ws1.Range("F13").Value = 403
$ws1.Range("F16").Value = 1373
$ws1.Range("F18").Value = 1640
$ws1.Range("F19").Value = 22
$ws1.Range("F23").Value = 48
$ws1.Range("F30").Value = 39
$ws1.Range("F32").Value = 3973
$ws1.Range("F36").Value = 1383
$ws1.Range("F38").Value = 1883

# Sheet "全部类型" (All types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1865
$ws4.Range("F14").Value = 403
$ws4.Range("F17").Value = 1373
$ws4.Range("F19").Value = 1640
$ws4.Range("F20").Value = 22
$ws4.Range("F24").Value = 48
$ws4.Range("F31").Value = 39
$ws4.Range("F33").Value = 3973
$ws4.Range("F39").Value = 1383
$ws4.Range("F41").Value = 1883
